# Add a new results row (row 5) to the "Resultados" sheet, mirroring the
# existing rows: torneo / categoria / grupo / posicion / jugador / dia_1 /
# dia_2 / total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"
$ws.Range("B5").Value = "Juveniles"
$ws.Range("C5").Value = "damas"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "Martínez, Valentina"
$ws.Range("F5").Value = 84
$ws.Range("H5").Value = 84
